$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# --- Row 3 header formatting -------------------------------------------------
# Make row 3 taller to fit the wrapped description text.
$ws.Rows.Item(3).RowHeight = 78.75

# G3 ("DESCRIPCIÓN" label) gets a left-only thin border plus wrap text.
$g3 = $ws.Range("G3")
$g3.Borders.LineStyle = -4142      # xlLineStyleNone - clear any existing border first
$g3.Borders.Item(7).LineStyle = 1  # xlEdgeLeft / xlContinuous - thin left border only
$g3.WrapText = $true

# H3 / I3 (rest of the merged G3:I3 region) just get wrap text enabled.
$ws.Range("H3").WrapText = $true
$ws.Range("I3").WrapText = $true

# --- Date value updates (August refresh) -------------------------------------
$ws.Range("B8").Value = 44652
$ws.Range("C8").Value = 44742
$ws.Range("AB8").Value = 44753
$ws.Range("AC8").Value = 44753

# --- Selection moves to C11 ---------------------------------------------------
$null = $ws.Range("C11").Select()

Write-Host "edit applied"
